$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the level layout with 4 more rows (20-24) -----------------
# Rows 20-23: repeating "_" wall-gap tiles in columns C, G, K, O
# (new shared string "PT" must be registered before "_" so the shared
#  string table order matches: index 10 = "PT", index 11 = "_")

$ws.Range("B24").Value = "PT"
$ws.Range("C24").Value = "PT"
$ws.Range("D24").Value = "PT"
$ws.Range("F24").Value = "PT"
$ws.Range("G24").Value = "PT"
$ws.Range("H24").Value = "PT"
$ws.Range("J24").Value = "PT"
$ws.Range("K24").Value = "PT"
$ws.Range("L24").Value = "PT"
$ws.Range("N24").Value = "PT"
$ws.Range("O24").Value = "PT"
$ws.Range("P24").Value = "PT"

$ws.Range("C20").Value = "_"
$ws.Range("G20").Value = "_"
$ws.Range("K20").Value = "_"
$ws.Range("O20").Value = "_"

$ws.Range("C21").Value = "_"
$ws.Range("G21").Value = "_"
$ws.Range("K21").Value = "_"
$ws.Range("O21").Value = "_"

$ws.Range("C22").Value = "_"
$ws.Range("G22").Value = "_"
$ws.Range("K22").Value = "_"
$ws.Range("O22").Value = "_"

$ws.Range("C23").Value = "_"
$ws.Range("G23").Value = "_"
$ws.Range("K23").Value = "_"
$ws.Range("O23").Value = "_"

# --- View/selection state ----------------------------------------------
$ws.Range("S21").Select()
